# Insert two new rows before row 93 (shifts old rows 93-128 down to 95-130)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("93:94").Insert()

# Common column values shared by every data row in this block
$commonA = 10
$commonB = "Vega Modelo de Temuco"
$commonC = "La Araucanía"
$commonE = 9
$commonF = "Fruta"
$commonG = 100103
$commonH = "Frutos de hueso (carozo)"
$commonI = 100103001
$commonJ = "Cereza"

# --- New row 93 ---
$ws.Cells.Item(93, 1).Value = $commonA
$ws.Cells.Item(93, 2).Value = $commonB
$ws.Cells.Item(93, 3).Value = $commonC
$ws.Cells.Item(93, 4).Value = 44524
$ws.Cells.Item(93, 5).Value = $commonE
$ws.Cells.Item(93, 6).Value = $commonF
$ws.Cells.Item(93, 7).Value = $commonG
$ws.Cells.Item(93, 8).Value = $commonH
$ws.Cells.Item(93, 9).Value = $commonI
$ws.Cells.Item(93, 10).Value = $commonJ
$ws.Cells.Item(93, 11).Value = "Early Burlat"
$ws.Cells.Item(93, 12).Value = "Primera"
$ws.Cells.Item(93, 13).Value = 365
$ws.Cells.Item(93, 14).Value = 20000
$ws.Cells.Item(93, 15).Value = 22000
$ws.Cells.Item(93, 16).Value = 21151
$ws.Cells.Item(93, 17).Value = "$/caja 10 kilos"
$ws.Cells.Item(93, 18).Value = "Región del Maule"
$ws.Cells.Item(93, 19).Value = 2115
$ws.Cells.Item(93, 20).Value = 10

# --- New row 94 ---
$ws.Cells.Item(94, 1).Value = $commonA
$ws.Cells.Item(94, 2).Value = $commonB
$ws.Cells.Item(94, 3).Value = $commonC
$ws.Cells.Item(94, 4).Value = 44524
$ws.Cells.Item(94, 5).Value = $commonE
$ws.Cells.Item(94, 6).Value = $commonF
$ws.Cells.Item(94, 7).Value = $commonG
$ws.Cells.Item(94, 8).Value = $commonH
$ws.Cells.Item(94, 9).Value = $commonI
$ws.Cells.Item(94, 10).Value = $commonJ
$ws.Cells.Item(94, 11).Value = "Rainier"
$ws.Cells.Item(94, 12).Value = "Primera"
$ws.Cells.Item(94, 13).Value = 110
$ws.Cells.Item(94, 14).Value = 40000
$ws.Cells.Item(94, 15).Value = 40000
$ws.Cells.Item(94, 16).Value = 40000
$ws.Cells.Item(94, 17).Value = "$/caja 12 kilos"
$ws.Cells.Item(94, 18).Value = "Región del Maule"
$ws.Cells.Item(94, 19).Value = 3333
$ws.Cells.Item(94, 20).Value = 12
